$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_5_2_0"
$ws.Range("B2").Value = 0.4991499343278779
$ws.Range("C2").Value = 0.797005131457062
$ws.Range("D2").Value = 0.9357095115723514
$ws.Range("E2").Value = 0.9098925090729403
$ws.Range("F2").Value = 0.5542934536933899
$ws.Range("G2").Value = 0.1713772863149643
$ws.Range("H2").Value = 0.09664588421583176
$ws.Range("I2").Value = 0.1362096071243286

$ws.Range("A3").Value = "model_5_2_1"
$ws.Range("B3").Value = 0.5178212898211507
$ws.Range("C3").Value = 0.7937577560265511
$ws.Range("D3").Value = 0.9331054203881135
$ws.Range("E3").Value = 0.9077139313395471
$ws.Range("F3").Value = 0.5336297750473022
$ws.Range("G3").Value = 0.1741188615560532
$ws.Range("H3").Value = 0.1005605310201645
$ws.Range("I3").Value = 0.1395028084516525

$ws.Range("A4").Value = "model_5_2_2"
$ws.Range("B4").Value = 0.5337416299968523
$ws.Range("C4").Value = 0.7860456471418314
$ws.Range("D4").Value = 0.9276256782757901
$ws.Range("E4").Value = 0.9028693633892489
$ws.Range("F4").Value = 0.516010582447052
$ws.Range("G4").Value = 0.1806297898292542
$ws.Range("H4").Value = 0.1087980568408966
$ws.Range("I4").Value = 0.1468260288238525

$ws.Range("A5").Value = "model_5_2_3"
$ws.Range("B5").Value = 0.5488227542589565
$ws.Range("C5").Value = 0.7774477560253389
$ws.Range("D5").Value = 0.9244183781006869
$ws.Range("E5").Value = 0.8988252579480505
$ws.Range("F5").Value = 0.4993202686309814
$ws.Range("G5").Value = 0.1878884881734848
$ws.Range("H5").Value = 0.1136194914579391
$ws.Range("I5").Value = 0.1529392451047897

$ws.Range("A6").Value = "model_5_2_4"
$ws.Range("B6").Value = 0.5622085774363264
$ws.Range("C6").Value = 0.7649178027755622
$ws.Range("D6").Value = 0.917052963310525
$ws.Range("E6").Value = 0.8916742955123218
$ws.Range("F6").Value = 0.4845061302185059
$ws.Range("G6").Value = 0.1984668374061584
$ws.Range("H6").Value = 0.1246916875243187
$ws.Range("I6").Value = 0.1637488603591919

$ws.Range("A7").Value = "model_5_2_5"
$ws.Range("B7").Value = 0.5750582138444471
$ws.Range("C7").Value = 0.7509115834815717
$ws.Range("D7").Value = 0.9128310659669285
$ws.Range("E7").Value = 0.8855572272658466
$ws.Range("F7").Value = 0.4702853858470917
$ws.Range("G7").Value = 0.2102915197610855
$ws.Range("H7").Value = 0.1310383379459381
$ws.Range("I7").Value = 0.1729956716299057

$ws.Range("A8").Value = "model_5_2_6"
$ws.Range("B8").Value = 0.6088663440977484
$ws.Range("C8").Value = 0.700899973931143
$ws.Range("D8").Value = 0.8992571866614022
$ws.Range("E8").Value = 0.8644174452489104
$ws.Range("F8").Value = 0.4328697025775909
$ws.Range("G8").Value = 0.2525135278701782
$ws.Range("H8").Value = 0.1514435261487961
$ws.Range("I8").Value = 0.204951286315918

$ws.Range("A9").Value = "model_5_2_7"
$ws.Range("B9").Value = 0.6191532532063224
$ws.Range("C9").Value = 0.6799147706495713
$ws.Range("D9").Value = 0.8939311660755833
$ws.Range("E9").Value = 0.8557200228616825
$ws.Range("F9").Value = 0.4214851260185242
$ws.Range("G9").Value = 0.2702301740646362
$ws.Range("H9").Value = 0.159449964761734
$ws.Range("I9").Value = 0.2180985957384109

$ws.Range("A10").Value = "model_5_2_8"
$ws.Range("B10").Value = 0.6289451754033903
$ws.Range("C10").Value = 0.6616958747432553
$ws.Range("D10").Value = 0.8885894509962513
$ws.Range("E10").Value = 0.8478330464993649
$ws.Range("F10").Value = 0.410648375749588
$ws.Range("G10").Value = 0.2856113612651825
$ws.Range("H10").Value = 0.1674799919128418
$ws.Range("I10").Value = 0.2300208359956741

$ws.Range("A11").Value = "model_5_2_9"
$ws.Range("B11").Value = 0.6372357319807127
$ws.Range("C11").Value = 0.6373264614484152
$ws.Range("D11").Value = 0.8800418984221485
$ws.Range("E11").Value = 0.8366275145055481
$ws.Range("F11").Value = 0.4014731049537659
$ws.Range("G11").Value = 0.3061851263046265
$ws.Range("H11").Value = 0.1803292632102966
$ws.Range("I11").Value = 0.2469594925642014

$ws.Range("A12").Value = "model_5_2_10"
$ws.Range("B12").Value = 0.6465544961375331
$ws.Range("C12").Value = 0.6076505061743377
$ws.Range("D12").Value = 0.867932426035404
$ws.Range("E12").Value = 0.8221861034424875
$ws.Range("F12").Value = 0.3911600112915039
$ws.Range("G12").Value = 0.331238865852356
$ws.Range("H12").Value = 0.1985330581665039
$ws.Range("I12").Value = 0.2687896192073822

$ws.Range("A13").Value = "model_5_2_11"
$ws.Range("B13").Value = 0.6551814428890996
$ws.Range("C13").Value = 0.5931200972065234
$ws.Range("D13").Value = 0.8624971885499081
$ws.Range("E13").Value = 0.815346314290913
$ws.Range("F13").Value = 0.3816125690937042
$ws.Range("G13").Value = 0.3435060977935791
$ws.Range("H13").Value = 0.2067036926746368
$ws.Range("I13").Value = 0.2791289091110229

$ws.Range("A14").Value = "model_5_2_12"
$ws.Range("B14").Value = 0.661735935290281
$ws.Range("C14").Value = 0.574306628325205
$ws.Range("D14").Value = 0.8539052255516685
$ws.Range("E14").Value = 0.8057630921668743
$ws.Range("F14").Value = 0.3743586242198944
$ws.Range("G14").Value = 0.3593892753124237
$ws.Range("H14").Value = 0.2196197062730789
$ws.Range("I14").Value = 0.2936152517795563

$ws.Range("A15").Value = "model_5_2_13"
$ws.Range("B15").Value = 0.6692032622201899
$ws.Range("C15").Value = 0.5599281123825547
$ws.Range("D15").Value = 0.8478726150355698
$ws.Range("E15").Value = 0.798688096261066
$ws.Range("F15").Value = 0.366094559431076
$ws.Range("G15").Value = 0.3715282678604126
$ws.Range("H15").Value = 0.2286883592605591
$ws.Range("I15").Value = 0.3043100833892822

$ws.Range("A16").Value = "model_5_2_14"
$ws.Range("B16").Value = 0.6764394974631465
$ws.Range("C16").Value = 0.5458626798353738
$ws.Range("D16").Value = 0.8425946874958261
$ws.Range("E16").Value = 0.7920600233237433
$ws.Range("F16").Value = 0.3580861389636993
$ws.Range("G16").Value = 0.3834029138088226
$ws.Range("H16").Value = 0.2366224825382233
$ws.Range("I16").Value = 0.3143292963504791

$ws.Range("A17").Value = "model_5_2_15"
$ws.Range("B17").Value = 0.6831406092900443
$ws.Range("C17").Value = 0.5307153039017063
$ws.Range("D17").Value = 0.8368854193734198
$ws.Range("E17").Value = 0.7849093409791503
$ws.Range("F17").Value = 0.3506700098514557
$ws.Range("G17").Value = 0.3961910009384155
$ws.Range("H17").Value = 0.2452050596475601
$ws.Range("I17").Value = 0.3251384794712067

$ws.Range("A18").Value = "model_5_2_16"
$ws.Range("B18").Value = 0.6898243966625583
$ws.Range("C18").Value = 0.5157521822073151
$ws.Range("D18").Value = 0.8315369589109818
$ws.Range("E18").Value = 0.7779816057238651
$ws.Range("F18").Value = 0.3432729840278625
$ws.Range("G18").Value = 0.408823549747467
$ws.Range("H18").Value = 0.2532452344894409
$ws.Range("I18").Value = 0.3356106579303741

$ws.Range("A19").Value = "model_5_2_17"
$ws.Range("B19").Value = 0.6959900100521925
$ws.Range("C19").Value = 0.5006328177481074
$ws.Range("D19").Value = 0.8269468693806176
$ws.Range("E19").Value = 0.7713631739279082
$ws.Range("F19").Value = 0.3364494740962982
$ws.Range("G19").Value = 0.4215879440307617
$ws.Range("H19").Value = 0.260145366191864
$ws.Range("I19").Value = 0.3456153869628906

$ws.Range("A20").Value = "model_5_2_18"
$ws.Range("B20").Value = 0.7011182497844446
$ws.Range("C20").Value = 0.4833495003981907
$ws.Range("D20").Value = 0.8208935719139676
$ws.Range("E20").Value = 0.7634197694521644
$ws.Range("F20").Value = 0.3307740390300751
$ws.Range("G20").Value = 0.4361792802810669
$ws.Range("H20").Value = 0.2692450881004333
$ws.Range("I20").Value = 0.3576228618621826

$ws.Range("A21").Value = "model_5_2_19"
$ws.Range("B21").Value = 0.7059720544492836
$ws.Range("C21").Value = 0.4658009814497532
$ws.Range("D21").Value = 0.8154456666350283
$ws.Range("E21").Value = 0.7556820923333235
$ws.Range("F21").Value = 0.3254022896289825
$ws.Range("G21").Value = 0.4509945511817932
$ws.Range("H21").Value = 0.277434766292572
$ws.Range("I21").Value = 0.3693193793296814

$ws.Range("A22").Value = "model_5_2_20"
$ws.Range("B22").Value = 0.7118482871427082
$ws.Range("C22").Value = 0.4572063423307237
$ws.Range("D22").Value = 0.8110951474436624
$ws.Range("E22").Value = 0.7511046741788743
$ws.Range("F22").Value = 0.318899005651474
$ws.Range("G22").Value = 0.4582505226135254
$ws.Range("H22").Value = 0.2839747667312622
$ws.Range("I22").Value = 0.3762387931346893

$ws.Range("A23").Value = "model_5_2_21"
$ws.Range("B23").Value = 0.7165270575771494
$ws.Range("C23").Value = 0.442927899471549
$ws.Range("D23").Value = 0.8064725871704768
$ws.Range("E23").Value = 0.7447196246840894
$ws.Range("F23").Value = 0.3137210011482239
$ws.Range("G23").Value = 0.4703050255775452
$ws.Range("H23").Value = 0.2909237146377563
$ws.Range("I23").Value = 0.3858906626701355

$ws.Range("A24").Value = "model_5_2_22"
$ws.Range("B24").Value = 0.7207970040506883
$ws.Range("C24").Value = 0.427482650221722
$ws.Range("D24").Value = 0.8016065161270935
$ws.Range("E24").Value = 0.7378754025132591
$ws.Range("F24").Value = 0.3089954555034637
$ws.Range("G24").Value = 0.4833446145057678
$ws.Range("H24").Value = 0.2982387244701385
$ws.Range("I24").Value = 0.3962366580963135

$ws.Range("A25").Value = "model_5_2_23"
$ws.Range("B25").Value = 0.7253019911840048
$ws.Range("C25").Value = 0.4154215902552625
$ws.Range("D25").Value = 0.7973530556015965
$ws.Range("E25").Value = 0.7323183953290837
$ws.Range("F25").Value = 0.3040097653865814
$ws.Range("G25").Value = 0.4935270547866821
$ws.Range("H25").Value = 0.304632842540741
$ws.Range("I25").Value = 0.4046368002891541

$ws.Range("A26").Value = "model_5_2_24"
$ws.Range("B26").Value = 0.7288922873632264
$ws.Range("C26").Value = 0.3982287815679891
$ws.Range("D26").Value = 0.7931465355561511
$ws.Range("E26").Value = 0.7252678475529202
$ws.Range("F26").Value = 0.3000363707542419
$ws.Range("G26").Value = 0.5080419778823853
$ws.Range("H26").Value = 0.3109563589096069
$ws.Range("I26").Value = 0.4152946472167969
